$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, matching the style of the existing header row (copy from F1)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Add the new data values
$ws.Range("G2").Value = 0.1256850772835605
$ws.Range("H2").Value = 0.99
